# Update "想去人数" (F column) figures for the 苏州-漫展信息 workbook.
# These updates mirror a re-generation of the source data (gh-pages build),
# bumping a handful of "want to go" counts on both the "展览" and
# "全部类型" sheets (which carry duplicate rows of the same events).

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览"     = @{
        3  = 599
        7  = 14966
        10 = 676
        11 = 15226
        12 = 36
        13 = 8740
        19 = 183
        20 = 12
        21 = 20
        22 = 510
        26 = 1085
        28 = 13
        34 = 29
        36 = 271
        37 = 426
        39 = 5368
    }
    "全部类型" = @{
        3  = 599
        7  = 14966
        10 = 676
        11 = 15226
        12 = 36
        13 = 8740
        20 = 183
        21 = 12
        22 = 20
        23 = 510
        27 = 1085
        29 = 13
        37 = 29
        39 = 271
        40 = 426
        42 = 5369
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
